$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Unprotect("")
$ws.Range("J12").Value = "Una máquina de escribir antigua y un computador portátil"
$ws.Rows.Item(12).AutoFit()
Write-Host "J12 set"
